# Update Shane Watson sheet: remove the "Oct 29 2020" match row, keeping
# only the "Oct 7 2020" match (previously row 3) now as row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite row 2 with the data that used to live in row 3. The numeric
# looking columns (G:K) must stay stored as text, so force text format
# before writing the values.
$ws.Range("G2:K2").NumberFormat = "@"

$ws.Range("A2").Value = " Oct 7 2020"
$ws.Range("B2").Value = " Abu Dhabi"
$ws.Range("C2").Value = "KKR won by 10 runs"
$ws.Range("D2").Value = "Kolkata Knight Riders"
$ws.Range("E2").Value = "Chennai Super Kings"
$ws.Range("F2").Value = "Shane Watson "
$ws.Range("G2").Value = "50"
$ws.Range("H2").Value = "40"
$ws.Range("I2").Value = "6"
$ws.Range("J2").Value = "1"
$ws.Range("K2").Value = "125.00"

# Remove the now-duplicate old row 3, shifting any rows below it upward
# (there are none, so this simply shrinks the used range to A1:K2).
$ws.Rows(3).Delete()
